# Updated cryptos list on Sun Jun  9 16:39:17 UTC 2024 with GitHub Actions
# Applies refreshed price / 1h-volume data (and a couple of row
# re-ranks: EthereumClassic <-> Binance-PegBSC-USD, Mantle <-> Monero)
# to the Sheet1 crypto table, cell-by-cell, matching the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.603.86'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '3.698.64'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5"); $c.Value = "'677.27"; $c.Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '
$c = $ws.Range("D6"); $c.Value = "'161.22"; $c.Style = "Normal"
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").Value = '  -0.07%  '
$c = $ws.Range("D8"); $c.Value = "'0.496"; $c.Style = "Normal"
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("E9").Value = '  +1.31%  '
$c = $ws.Range("D10"); $c.Value = "'7.11"; $c.Style = "Normal"
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  +0.63%  '
$c = $ws.Range("D13"); $c.Value = "'32.57"; $c.Style = "Normal"
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").Value = '3.702.47'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '69.569.12'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("E16").Value = '  +2.37%  '
$c = $ws.Range("D17"); $c.Value = "'16.03"; $c.Style = "Normal"
$ws.Range("E17").Value = '  +1.07%  '
$c = $ws.Range("D18"); $c.Value = "'6.49"; $c.Style = "Normal"
$ws.Range("E18").Value = '  +0.74%  '
$c = $ws.Range("D19"); $c.Value = "'471.96"; $c.Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$c = $ws.Range("D20"); $c.Value = "'9.88"; $c.Style = "Normal"
$ws.Range("E20").Value = '  -2.02%  '
$c = $ws.Range("D21"); $c.Value = "'0.651"; $c.Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '
$c = $ws.Range("D22"); $c.Value = "'80.55"; $c.Style = "Normal"
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("D23").Value = '3.842.81'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +1.98%  '
$c = $ws.Range("D26"); $c.Value = "'10.89"; $c.Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$c = $ws.Range("D27"); $c.Value = "'9.14"; $c.Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '
$c = $ws.Range("D28"); $c.Value = "'2.71"; $c.Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("E30").Value = '  +0.74%  '
$c = $ws.Range("D31"); $c.Value = "'6.61"; $c.Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D32"); $c.Value = "'27.02"; $c.Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D33"); $c.Value = "'1.00"; $c.Style = "Normal"
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '3.686.44'
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("E35").Value = '  +0.62%  '
$c = $ws.Range("D36"); $c.Value = "'8.49"; $c.Style = "Normal"
$ws.Range("E36").Value = '  +4.21%  '
$ws.Range("E37").Value = '  +1.56%  '
$c = $ws.Range("D39"); $c.Value = "'2.26"; $c.Style = "Normal"
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("E40").Value = '  -0.07%  '
$c = $ws.Range("D41"); $c.Value = "'0.0903"; $c.Style = "Normal"
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D42"); $c.Value = "'0.944"; $c.Style = "Normal"
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D43"); $c.Value = "'166.99"; $c.Style = "Normal"
$ws.Range("E43").Value = '  +1.21%  '
$c = $ws.Range("D44"); $c.Value = "'46.41"; $c.Style = "Normal"
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("E45").Value = '  +1.57%  '
$c = $ws.Range("D46"); $c.Value = "'0.000279"; $c.Style = "Normal"
$ws.Range("E46").Value = '  +2.16%  '
$c = $ws.Range("D47"); $c.Value = "'28.14"; $c.Style = "Normal"
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  -1.28%  '
$c = $ws.Range("D50"); $c.Value = "'7.89"; $c.Style = "Normal"
$ws.Range("E51").Value = '  +2.10%  '
